$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.972.79'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.504.30'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.78%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '535.07'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +5.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.25'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.568'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.515.87'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.78%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.59%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.19'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.332'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.961.54'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.932.03'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.44'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.06%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.514.73'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.66'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.27'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '321.08'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.27'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +8.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.03'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.411'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.00%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.53'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.46%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +5.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '172.16'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.27%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.19'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.71%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.33%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.18%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.84%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.99'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.80%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.69%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'SuiNetwork'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.819'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +8.98%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.51'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.49'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.80%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '277.60'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.66%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '132.26'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +10.94%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.07'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.82%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.41%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0513'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.79%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.90'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.752.10'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.29%  '
